$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Row, $Col, $Text)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-CellText 2 4 '29.436.43'
Set-CellText 2 5 '  +0.11%  '
Set-CellText 3 4 '1.850.39'
Set-CellText 3 5 '  +0.12%  '
Set-CellText 4 4 '0.9992'
Set-CellText 4 5 '  -0.04%  '
Set-CellText 5 4 '240.27'
Set-CellText 5 5 '  -0.12%  '
Set-CellText 6 4 '0.6294'
Set-CellText 6 5 '  -0.25%  '
Set-CellText 7 4 '1.000'
Set-CellText 7 5 '  +0.01%  '
Set-CellText 8 4 '0.07658'
Set-CellText 8 5 '  +0.67%  '
Set-CellText 9 4 '0.2923'
Set-CellText 9 5 '  -0.55%  '
Set-CellText 10 4 '24.79'
Set-CellText 10 5 '  +0.77%  '
Set-CellText 11 4 '0.07753'
Set-CellText 11 5 '  +0.11%  '
Set-CellText 12 4 '1.917.25'
Set-CellText 12 5 '  +1.17%  '
Set-CellText 13 4 '5.036'
Set-CellText 13 5 '  +0.53%  '
Set-CellText 14 4 '0.6817'
Set-CellText 14 5 '  +0.10%  '
Set-CellText 15 4 '0.00001052'
Set-CellText 15 5 '  -3.43%  '
Set-CellText 16 4 '83.61'
Set-CellText 16 5 '  -0.08%  '
Set-CellText 17 2 'Uniswap'
Set-CellText 17 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText 17 4 '6.194'
Set-CellText 17 5 '  -0.07%  '
Set-CellText 18 2 'WrappedBTC'
Set-CellText 18 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText 18 4 '29.432.70'
Set-CellText 18 5 '  -0.04%  '
Set-CellText 19 2 'BitcoinCash'
Set-CellText 19 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText 19 4 '229.53'
Set-CellText 19 5 '  -0.11%  '
Set-CellText 20 2 'Avalanche'
Set-CellText 20 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-CellText 20 4 '12.35'
Set-CellText 20 5 '  -0.79%  '
Set-CellText 21 2 'Dai'
Set-CellText 21 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText 21 4 '1.000'
Set-CellText 21 5 '  -0.01%  '
Set-CellText 22 2 'Chainlink'
Set-CellText 22 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText 22 4 '7.491'
Set-CellText 22 5 '  -0.66%  '
Set-CellText 23 2 'BinanceUSD'
Set-CellText 23 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-CellText 23 4 '1.000'
Set-CellText 23 5 '  -0.04%  '
Set-CellText 24 2 'Monero'
Set-CellText 24 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText 24 4 '157.29'
Set-CellText 24 5 '  +0.00%  '
Set-CellText 25 2 'Stellar'
Set-CellText 25 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText 25 4 '0.1385'
Set-CellText 25 5 '  -1.26%  '
Set-CellText 26 2 'Cosmos'
Set-CellText 26 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText 26 4 '8.422'
Set-CellText 26 5 '  +0.69%  '
Set-CellText 27 2 'EthereumClassic'
Set-CellText 27 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-CellText 27 4 '17.79'
Set-CellText 27 5 '  +0.77%  '
Set-CellText 28 2 'Toncoin'
Set-CellText 28 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText 28 4 '1.363'
Set-CellText 28 5 '  +4.88%  '
Set-CellText 29 2 'PancakeSwap'
Set-CellText 29 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText 29 4 '1.462'
Set-CellText 29 5 '  -0.21%  '
Set-CellText 30 2 'Hedera'
Set-CellText 30 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText 30 4 '0.05612'
Set-CellText 30 5 '  +0.26%  '
Set-CellText 31 2 'Filecoin'
Set-CellText 31 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText 31 4 '4.131'
Set-CellText 31 5 '  +0.32%  '
Set-CellText 32 2 'InternetComputer(DFINITY)'
Set-CellText 32 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText 32 4 '4.060'
Set-CellText 32 5 '  +0.47%  '
Set-CellText 33 2 'LidoDAOToken'
Set-CellText 33 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText 33 4 '1.844'
Set-CellText 33 5 '  -0.84%  '
Set-CellText 34 2 'ARBITRUM'
Set-CellText 34 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText 34 4 '1.167'
Set-CellText 34 5 '  +0.73%  '
Set-CellText 35 2 'ImmutableX'
Set-CellText 35 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText 35 4 '0.7004'
Set-CellText 35 5 '  -1.81%  '
Set-CellText 36 2 'HuobiToken'
Set-CellText 36 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText 36 4 '2.593'
Set-CellText 36 5 '  +0.14%  '
Set-CellText 37 2 'VeChain'
Set-CellText 37 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText 37 4 '0.01802'
Set-CellText 37 5 '  -0.46%  '
Set-CellText 38 4 '1.223.34'
Set-CellText 38 5 '  -1.84%  '
Set-CellText 39 2 'MXToken'
Set-CellText 39 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText 39 4 '2.746'
Set-CellText 39 5 '  -1.08%  '
Set-CellText 40 2 'FraxShare'
Set-CellText 40 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText 40 4 '6.445'
Set-CellText 40 5 '  +0.17%  '
Set-CellText 41 2 'TrustWalletToken'
Set-CellText 41 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText 41 4 '0.9075'
Set-CellText 41 5 '  +0.56%  '
Set-CellText 42 2 'PaxDollar'
Set-CellText 42 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-CellText 42 4 '1.000'
Set-CellText 42 5 '  +0.00%  '
Set-CellText 43 2 'Quant'
Set-CellText 43 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText 43 4 '101.95'
Set-CellText 43 5 '  +0.02%  '
Set-CellText 44 2 'Aave'
Set-CellText 44 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText 44 4 '66.34'
Set-CellText 44 5 '  +0.51%  '
Set-CellText 45 2 'Aptos'
Set-CellText 45 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText 45 4 '7.204'
Set-CellText 45 5 '  +0.67%  '
Set-CellText 46 2 'BabyDogeCoin'
Set-CellText 46 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText 46 4 '0.00000000120'
Set-CellText 46 5 '  +0.07%  '
Set-CellText 47 2 'TheSandbox'
Set-CellText 47 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText 47 4 '0.4025'
Set-CellText 47 5 '  +0.08%  '
Set-CellText 48 2 'EnergySwap'
Set-CellText 48 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 48 4 '9.055'
Set-CellText 48 5 '  +0.87%  '
Set-CellText 49 2 'Algorand'
Set-CellText 49 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText 49 4 '0.1157'
Set-CellText 49 5 '  +3.07%  '
Set-CellText 50 2 'RenderToken'
Set-CellText 50 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText 50 4 '1.678'
Set-CellText 50 5 '  -0.80%  '
Set-CellText 51 2 'Cronos'
Set-CellText 51 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText 51 4 '0.05709'
Set-CellText 51 5 '  -0.07%  '
